# Add the new ADC power-down (PD) register row, and move the existing
# "clock select" register row from adr=16 (row 23) down to adr=32 (row 39).
#
# commit message: "add ADC pd register"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Move the "clock select" register (currently documented at row 23,
#    register address 16) down to row 39 (register address 32) -- the
#    firmware moved that register's address.
$ws.Range("B39").Value = $ws.Range("B23").Value2
$ws.Range("C39").Value = $ws.Range("C23").Value2
$ws.Range("D39").Value = $ws.Range("D23").Value2
$ws.Range("E39").Value = $ws.Range("E23").Value2

# 2. Clear the old location. B23/C23 and E23 become completely empty
#    (no leftover formatting), while D23 keeps its original formatting
#    but with no value.
$ws.Range("B23:C23").Clear()
$ws.Range("E23").Clear()
$ws.Range("D23").ClearContents()

# 3. Document the new "adc->powerdown (pd)" register at row 67
#    (register address 60).
$ws.Range("B67").Value = "adc->powerdown (pd)"
$ws.Range("C67").Value = "lower 4 bits set PD (I+Q channel, ADCs 0->3)"
$ws.Range("D67").Value = "0x000000"
